$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'72.538.99"
$ws.Range("E2").Value = "'  +1.82%  "

$ws.Range("D3").Value = "'2.690.07"
$ws.Range("E3").Value = "'  +2.62%  "

$ws.Range("E4").Value = "'  -0.03%  "

$ws.Range("D5").Value = "'601.65"
$ws.Range("E5").Value = "'  -0.75%  "

$ws.Range("D6").Value = "'178.04"
$ws.Range("E6").Value = "'  -1.60%  "

$ws.Range("E7").Value = "'  -0.02%  "

$ws.Range("E8").Value = "'  +0.20%  "

$ws.Range("D9").Value = "'2.685.94"
$ws.Range("E9").Value = "'  +2.48%  "

$ws.Range("D10").Value = "'0.170"
$ws.Range("E10").Value = "'  +3.88%  "

$ws.Range("E12").Value = "'  +2.70%  "

$ws.Range("D13").Value = "'5.06"
$ws.Range("E13").Value = "'  +0.86%  "

$ws.Range("D14").Value = "'3.166.72"
$ws.Range("E14").Value = "'  +2.82%  "

$ws.Range("E15").Value = "'  +0.30%  "

$ws.Range("D16").Value = "'72.457.78"
$ws.Range("E16").Value = "'  +1.76%  "

$ws.Range("D17").Value = "'26.43"
$ws.Range("E17").Value = "'  -0.52%  "

$ws.Range("D18").Value = "'2.678.30"
$ws.Range("E18").Value = "'  +1.51%  "

$ws.Range("D19").Value = "'12.00"
$ws.Range("E19").Value = "'  +4.57%  "

$ws.Range("D20").Value = "'7.95"
$ws.Range("E20").Value = "'  +0.63%  "

$ws.Range("D21").Value = "'372.80"
$ws.Range("E21").Value = "'  -2.96%  "

$ws.Range("D22").Value = "'4.18"
$ws.Range("E22").Value = "'  +0.98%  "

$ws.Range("E23").Value = "'  +8.60%  "

$ws.Range("E24").Value = "'  +0.20%  "

$ws.Range("E25").Value = "'  -0.04%  "

$ws.Range("D26").Value = "'4.35"
$ws.Range("E26").Value = "'  -2.57%  "

$ws.Range("D27").Value = "'9.89"
$ws.Range("E27").Value = "'  +2.37%  "

$ws.Range("D28").Value = "'2.823.25"
$ws.Range("E28").Value = "'  +2.38%  "

$ws.Range("E29").Value = "'  +0.12%  "

$ws.Range("D30").Value = "'0.0₃0944"
$ws.Range("E30").Value = "'  -2.02%  "

$ws.Range("D31").Value = "'8.13"
$ws.Range("E31").Value = "'  +0.95%  "

$ws.Range("D32").Value = "'515.97"
$ws.Range("E32").Value = "'  -5.14%  "

$ws.Range("D33").Value = "'1.31"
$ws.Range("E33").Value = "'  -0.88%  "

$ws.Range("E34").Value = "'  -0.63%  "

$ws.Range("E35").Value = "'  -0.11%  "

$ws.Range("D36").Value = "'162.84"
$ws.Range("E36").Value = "'  -1.81%  "

$ws.Range("D37").Value = "'19.63"
$ws.Range("E37").Value = "'  +2.31%  "

$ws.Range("E38").Value = "'  +0.54%  "

$ws.Range("E39").Value = "'  -0.21%  "

$ws.Range("D40").Value = "'1.82"
$ws.Range("E40").Value = "'  -3.26%  "

$ws.Range("D41").Value = "'0.108"
$ws.Range("E41").Value = "'  -8.22%  "

$ws.Range("E42").Value = "'  +0.04%  "

$ws.Range("E43").Value = "'  -0.14%  "

$ws.Range("E44").Value = "'  -2.35%  "

$ws.Range("E45").Value = "'  +0.93%  "

$ws.Range("D46").Value = "'39.28"
$ws.Range("E46").Value = "'  -1.81%  "

$ws.Range("D47").Value = "'154.22"
$ws.Range("E47").Value = "'  -0.08%  "

$ws.Range("D48").Value = "'3.74"
$ws.Range("E48").Value = "'  +2.89%  "

$ws.Range("D49").Value = "'0.554"
$ws.Range("E49").Value = "'  +3.88%  "

$ws.Range("D50").Value = "'1.74"
$ws.Range("E50").Value = "'  +2.30%  "

$ws.Range("D51").Value = "'0.0768"
$ws.Range("E51").Value = "'  +1.66%  "
